# A new September log entry ("beneficiary" at 2024-09-05 17:04:10) was added at
# the top of the running log on the "2024" sheet. That pushes every later
# September entry (columns R/S) down by one row, which in turn pushes the
# August entries (columns P/Q, rows 76-79) and the trailing "Broadband" group
# label (column A) down by one row as well, growing the sheet from 79 to 80
# rows.
#
# Rather than using Rows.Insert() (which, on this engine, only materialises
# cells that already carry a real value - every placeholder "blank" cell in
# the untouched columns of a brand-new row would silently vanish), we shift
# the data by writing the new value into every cell whose content actually
# changes, bottom constant columns stay exactly as they are. This keeps every
# already-existing row fully intact.
#
# Row 80 is genuinely new, so to keep its placeholder cells (columns that
# stay blank) present in the file, we materialise it by copying the
# still-original row 79 (A79="Broadband", everything else blank) down onto
# it with a single Range.Copy(destination) - that preserves blank cell
# records the same way real Excel would when you drag a row down. Only then
# do we overwrite row 79 with its own new contents.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# --- materialise the new row 80 from the current (pre-edit) row 79 ---------
$ws.Range("A79:Y79").Copy($ws.Range("A80"))

# --- pre-seed the cells that will become blank so their (blank) cell record
#     stays present, by copying existing blank cells onto them -------------
$ws.Range("B75:C75").Copy($ws.Range("P75"))
$ws.Range("B79").Copy($ws.Range("A79"))

# --- new top September entry -----------------------------------------------
$ws.Range("R31").Value = "beneficiary"
$ws.Range("S31").Value = "2024-09-05 17:04:10"

# --- September log (R/S) cascades down by one row, rows 32-75 --------------
$ws.Range("R32").Value = "bal axisbank"
$ws.Range("S32").Value = "2024-09-05 16:52:25"
$ws.Range("R33").Value = "share anyone axis"
$ws.Range("S33").Value = "2024-09-05 16:38:59"
$ws.Range("R34").Value = "transfer anyone axis"
$ws.Range("S34").Value = "2024-09-05 16:35:58"
$ws.Range("R35").Value = "share anyone axis"
$ws.Range("S35").Value = "2024-09-05 16:31:34"
$ws.Range("R36").Value = "transfer"
$ws.Range("S36").Value = "2024-09-05 16:28:38"
$ws.Range("R37").Value = "bal axisbank axis"
$ws.Range("S37").Value = "2024-09-05 16:26:56"
$ws.Range("R38").Value = "bal axisbank"
$ws.Range("S38").Value = "2024-09-05 16:26:55"
$ws.Range("S39").Value = "2024-09-05 16:25:07"
$ws.Range("R40").Value = "transfer"
$ws.Range("S40").Value = "2024-09-05 16:22:23"
$ws.Range("R41").Value = "share anyone axis"
$ws.Range("S41").Value = "2024-09-05 16:06:05"
$ws.Range("R42").Value = "internet bal axisbank"
$ws.Range("S42").Value = "2024-09-05 16:05:55"
$ws.Range("R43").Value = "transfer share anyone axis"
$ws.Range("S43").Value = "2024-09-05 16:03:14"
$ws.Range("R44").Value = "axis"
$ws.Range("R45").Value = "your net internet"
$ws.Range("S45").Value = "2024-09-05 15:57:15"
$ws.Range("R46").Value = "hear your feedback atm"
$ws.Range("S46").Value = "2024-09-05 14:21:08"
$ws.Range("S47").Value = "2024-09-05 14:18:32"
$ws.Range("S48").Value = "2024-09-05 14:13:16"
$ws.Range("R49").Value = "axis bna"
$ws.Range("S49").Value = "2024-09-05 14:15:23"
$ws.Range("R50").Value = "balance your axis"
$ws.Range("S50").Value = "2024-09-05 09:20:57"
$ws.Range("R51").Value = "bal axis"
$ws.Range("S51").Value = "2024-09-05 09:06:25"
$ws.Range("R52").Value = "broker"
$ws.Range("S52").Value = "2024-09-04 21:20:47"
$ws.Range("R53").Value = "exclusive on axis"
$ws.Range("S53").Value = "2024-09-04 13:21:05"
$ws.Range("R54").Value = "your corporate axis"
$ws.Range("S54").Value = "2024-09-04 11:46:10"
$ws.Range("R55").Value = "balance your axis"
$ws.Range("S55").Value = "2024-09-04 08:14:16"
$ws.Range("R56").Value = "axis"
$ws.Range("S56").Value = "2024-09-04 07:02:13"
$ws.Range("R57").Value = "bal axisbank w axis"
$ws.Range("S57").Value = "2024-09-04 06:53:15"
$ws.Range("R58").Value = "logging iob internet"
$ws.Range("S58").Value = "2024-09-03 20:09:12"
$ws.Range("R59").Value = "password internet"
$ws.Range("S59").Value = "2024-09-03 20:05:31"
$ws.Range("R60").Value = "logging iob internet"
$ws.Range("S60").Value = "2024-09-03 20:05:09"
$ws.Range("R61").Value = "internet"
$ws.Range("S61").Value = "2024-09-03 19:58:18"
$ws.Range("S62").Value = "2024-09-03 19:54:49"
$ws.Range("R63").Value = "login internet invalid"
$ws.Range("S63").Value = "2024-09-03 19:56:17"
$ws.Range("R64").Value = "corporate internet share"
$ws.Range("S64").Value = "2024-09-03 19:22:58"
$ws.Range("R65").Value = "login sbi internet personal do not share anyone"
$ws.Range("S65").Value = "2024-09-03 19:17:10"
$ws.Range("R66").Value = "login internet personal share"
$ws.Range("S66").Value = "2024-09-03 19:13:40"
$ws.Range("R67").Value = "internet verify it"
$ws.Range("S67").Value = "2024-09-03 19:05:49"
$ws.Range("R68").Value = "balance your axis"
$ws.Range("S68").Value = "2024-09-03 13:14:06"
$ws.Range("R69").Value = "lounge"
$ws.Range("S69").Value = "2024-09-03 13:08:08"
$ws.Range("R70").Value = "balance your axis"
$ws.Range("S70").Value = "2024-09-03 11:21:30"
$ws.Range("R71").Value = "broker"
$ws.Range("S71").Value = "2024-09-01 22:35:38"
$ws.Range("S72").Value = "2024-09-01 10:12:03"
$ws.Range("S73").Value = "2024-09-01 09:42:38"
$ws.Range("S74").Value = "2024-09-01 09:29:24"

# --- row 75: August (P/Q) cleared, gains the last September entry (R/S) ----
$ws.Range("P75").Value = ""
$ws.Range("Q75").Value = ""
$ws.Range("R75").Value = "amazeloan"
$ws.Range("S75").Value = "2024-09-01 09:27:06"

# --- August log (P/Q) cascades down by one row, rows 76-79 ------------------
$ws.Range("Q76").Value = "2024-08-30 12:15:48"
$ws.Range("Q77").Value = "2024-08-21 20:17:10"
$ws.Range("Q78").Value = "2024-08-21 20:16:45"
$ws.Range("P79").Value = "hdfc"
$ws.Range("Q79").Value = "2024-08-21 20:15:50"

# --- the "Broadband" group label moves from row 79 down to the new row 80 --
$ws.Range("A79").Value = ""
$ws.Range("A80").Value = "Broadband"
